$d = $word.ActiveDocument

# --- Change 1: rewrite the "Item based" intro paragraph (5 runs -> 1 run, new wording) ---
$oldP1 = 'This item based recommender first looks up what artists user U already follows. Those artists that U follows and their listen count are added to a list, sorted by listencount from high to low. Then using an item similarity pickle, all the similar artists to those that user U already follows are found. Then a top N is made from the most similar artists to the ones already followed. The count that is kept track of together with the similarity scores, define a weight for the similar artists. Then using these weights another top N is made for the similar artists that user U should follow but doesn’t follow yet.'
$newP1 = 'The most crucial step in the item based recommender is to compute the similarity scores using cosine similarity. The file containg these similarities between artists is called “item_sim50” and only contains similarity scores higher than 0.50 as we found out that none of the results we first had, using a file containg only scores higher than 0.10, didn’t contain any similarity scores lower than 0.50. Setting this value higher, made the similarity file much smaller, going from 23 million lines to 1.1 million lines. This in turn made running the program much faster. '
$found1 = $d.Content.Find.Execute($oldP1, $true, $false, $false, $false, $false, $true, 1, $false, $newP1, 2)
Write-Output "Change1 found/replaced: $found1"

# --- Change 2: split the "For the optimal value of N..." paragraph into three paragraphs ---
$targetPara = $d.Paragraphs(25)
Write-Output "Target paragraph text: $($targetPara.Range.Text)"
$r = $targetPara.Range.Duplicate
$r.Collapse(1)
$xmlPayload = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:ind w:hanging="0"/><w:rPr/></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">Our item based recommender first looks up what artists user U already follows, using the training dataset. These artists and their listencount are added to a list for the user U, sorted in a top N by listencount from highest to lowest. Then using the item similarity pickle explained above, similar artists to the ones U already follows could be </w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>computed. These also were sorted in a top N from highest to lowest on their similarity scores. Using the listen count for each already listened artist and the similarity scores for the similar artists, a weight is calculated by multiplying the listencount by the similarityscores. This is done to only get the most similar artists to artists that user U listens to a lot, but doesn’t know yet. This weighted score for each similar artist is sorted in a top 10 and this in turn is the end result for the recommendation of artists that user U doesn’t listen to yet, but might want to.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:ind w:hanging="0"/><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr/></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:ind w:hanging="0"/><w:rPr/></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">For the optimal value of N </w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>it was found that</w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> N=30 give</w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>s the</w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> most desirable results for speed and score.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xmlPayload) | Out-Null

Write-Output "Done; paragraph count now: $($d.Paragraphs.Count)"
